$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New query text blocks (instrument model DNBSEQ-G400 -> Illumina HiSeq 2500)
$statsQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@
$participantsQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@
$samplesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$filesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@

$neo4jFile = 'TC06_CDS_Filter_InstrumentModel-Illumina HiSeq 2500_Neo4jData.xlsx'
$webFile = 'TC06_CDS_Filter_InstrumentModel-Illumina HiSeq 2500_WebData.xlsx'

# Row 2 - ParticipantsTab
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statsQuery
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

# Row 3 - SamplesTab
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statsQuery
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

# Row 4 - FilesTab
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statsQuery
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# Column width adjustments for columns D and E (target stored widths: 94.85546875 / 93.140625;
# ColumnWidth is quantized internally, so use the nearest settable values)
$ws.Columns.Item(4).ColumnWidth = 94
$ws.Columns.Item(5).ColumnWidth = 92.333333

# Update selection to C2
$ws.Range("C2").Select()
